$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.548.81"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "2.624.42"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.83"
$ws.Range("E5").Value = "  +2.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.47"
$ws.Range("E6").Value = "  +1.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.55"
$ws.Range("E9").Value = "  -0.49%  "

$ws.Range("E10").Value = "  +0.59%  "

$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("E12").Value = "  +3.44%  "

$ws.Range("D13").Value = "3.091.63"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.13"
$ws.Range("E14").Value = "  +11.79%  "

$ws.Range("D15").Value = "60.545.91"
$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").Value = "2.630.37"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.54"
$ws.Range("E18").Value = "  +2.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.70"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.38"
$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("E21").Value = "  -1.59%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.531"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.88"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("E26").Value = "  +1.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.10"
$ws.Range("E27").Value = "  +5.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.03"
$ws.Range("E28").Value = "  +11.17%  "

$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  +1.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.65"
$ws.Range("E30").Value = "  +5.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.09"
$ws.Range("E31").Value = "  +5.05%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.61"
$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("E34").Value = "  +9.86%  "

$ws.Range("E35").Value = "  +4.73%  "

$ws.Range("E36").Value = "  +7.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.64"
$ws.Range("E37").Value = "  +2.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "331.38"
$ws.Range("E38").Value = "  +12.74%  "

$ws.Range("E39").Value = "  +4.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.36"
$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.875"
$ws.Range("E41").Value = "  +3.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.13"
$ws.Range("E42").Value = "  +6.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.70"
$ws.Range("E43").Value = "  +4.28%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0998"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "133.00"
$ws.Range("E45").Value = "  -3.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.08"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0556"
$ws.Range("E48").Value = "  +2.07%  "

$ws.Range("E49").Value = "  +0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0244"
$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("E51").Value = "  +0.18%  "
